$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Cracker"
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = "BenFred"

$ws.Range("F5").Select()
